# Update countries & provincias Spain
# Applies:
#  - Swap the "Kazajistan" / "Rumania" rows (Rumania's case count overtook
#    Kazajistan's, so it now sorts above it) at rows 33/34.
#  - Refresh numeric stats (Casos totales, Nuevos casos, Casos activos,
#    Recuperados, Casos criticos, Muertes) for several countries.
#  - Refresh the "Datos actualizados a ..." timestamp cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp banner (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 17 de Septiembre de 2020 a las 12:24"

# --- Country name swap: row 33 becomes Rumania, row 34 becomes Kazajistan ---
$ws.Range("A33").Value = "Rumania"
$ws.Range("A34").Value = "Kazajistan"

# --- Row 15 (Iran) ---
$ws.Range("B15").Value = 413149
$ws.Range("C15").Value = 2815
$ws.Range("D15").Value = 353848
$ws.Range("E15").Value = 35493
$ws.Range("G15").Value = 176
$ws.Range("H15").Value = 23808

# --- Row 18 (Banglades) ---
$ws.Range("B18").Value = 344264
$ws.Range("C18").Value = 1593
$ws.Range("D18").Value = 250412
$ws.Range("E18").Value = 88993
$ws.Range("G18").Value = 36
$ws.Range("H18").Value = 4859

# --- Row 31 (Catar) ---
$ws.Range("B31").Value = 122693
$ws.Range("C31").Value = 244
$ws.Range("D31").Value = 119613
$ws.Range("E31").Value = 2872

# --- Row 33 (now Rumania, after name swap above) ---
$ws.Range("B33").Value = 108690
$ws.Range("C33").Value = 1679
$ws.Range("D33").Value = 43244
$ws.Range("E33").Value = 61134
$ws.Range("G33").Value = 27
$ws.Range("H33").Value = 4312

# --- Row 34 (now Kazajistan, after name swap above) ---
$ws.Range("B34").Value = 107056
$ws.Range("C34").Value = 72
$ws.Range("D34").Value = 101455
$ws.Range("E34").Value = 3930
$ws.Range("H34").Value = 1671

# --- Row 70 (Austria) ---
$ws.Range("B70").Value = 35853
$ws.Range("C70").Value = 780
$ws.Range("D70").Value = 28044
$ws.Range("E70").Value = 7051

# --- Row 97 (Malasia) ---
$ws.Range("B97").Value = 10052
$ws.Range("C97").Value = 21
$ws.Range("D97").Value = 9250
$ws.Range("E97").Value = 674
